# "fixed no fault data" - refresh the simulation results on Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Updated Message Exchanged (E), Bytes Transferred (F) and Response Time (G)
# values for rows 2-10, replacing the earlier (faulty) measurements.
$data = @{
    2  = @{ E = 1142; F = 34688896; G = 11184648 }
    3  = @{ E = 1392; F = 34669104; G = 49798326 }
    4  = @{ E = 1561; F = 37041440; G = 92016487 }
    5  = @{ E = 1727; F = 38325240; G = 138734527 }
    6  = @{ E = 2072; F = 41786488; G = 184429954 }
    7  = @{ E = 2204; F = 42023400; G = 227660000 }
    8  = @{ E = 2393; F = 45870480; G = 275262240 }
    9  = @{ E = 2679; F = 47429264; G = 328379792 }
    10 = @{ E = 2918; F = 51596320; G = 369418149 }
}

foreach ($row in ($data.Keys | Sort-Object)) {
    $vals = $data[$row]
    $ws.Range("E$row").Value = $vals.E
    $ws.Range("F$row").Value = $vals.F
    $ws.Range("G$row").Value = $vals.G
}

# The wider values in column F (Bytes Transferred) now need the column
# to be auto-sized to fit its contents, just like column G already is.
$ws.Columns.Item(6).AutoFit() | Out-Null

# Move/collapse the selection from the old C2:C10 range down to G10.
$ws.Range("G10").Select() | Out-Null

# Keep the workbook window geometry in sync with the authored file.
try {
    $win = $wb.Windows.Item(1)
    $win.Left = 18380
    $win.Top = 1500
    $win.Width = 15220
    $win.Height = 17440
} catch {}

Write-Host "Applied no-fault data fix"
